# Apply edits described by the commit diff:
#  1. Rename header "Requested quantity" -> "Weekly_PO_Qty" on "Weekly Quantity" sheet
#  2. Rename header "Requested quantity" -> "Monthly_PO_Qty" on "Monthly Trend" sheet
#  3. Add a new "PO Forecast" worksheet (after "Monthly Trend") with forecast data

$wb = $excel.ActiveWorkbook

# --- 1. Update "Weekly Quantity" header ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Update "Monthly Trend" header ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add new "PO Forecast" worksheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy header formatting (bold, centered, bordered) from "Weekly Quantity" sheet
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy date-formatted style for column A from "Weekly Quantity" sheet
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A22").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(45361.99999999999, 9, -17.65803011345704, 34.32880897832849),
    @(45382.99999999999, 10, -15.82609942668298, 33.30947364770498),
    @(45459.99999999999, 14, -9.193175961991718, 39.25819574264543),
    @(45487.99999999999, 16, -8.929270805954582, 40.34797955842701),
    @(45515.99999999999, 17, -8.265940654728816, 43.07445160124044),
    @(45522.99999999999, 18, -8.024283550569264, 42.49166109349532),
    @(45529.99999999999, 18, -9.185830716550086, 44.00222930756258),
    @(45536.99999999999, 18, -7.656485084964252, 43.05267837482739),
    @(45543.99999999999, 19, -5.760830549030146, 43.44600421922447),
    @(45550.99999999999, 19, -6.739642272186273, 42.66724327756742),
    @(45557.99999999999, 19, -8.059742736594536, 44.9155749843389),
    @(45571.99999999999, 20, -4.592077043791006, 44.75856737829242),
    @(45578.99999999999, 21, -4.371231302353982, 44.35478614916112),
    @(45585.99999999999, 21, -3.319396230902302, 46.82654264066245),
    @(45592.99999999999, 21, -3.318305014230794, 44.40043885823146),
    @(45599.99999999999, 22, -2.475222516951801, 47.44608812930522),
    @(45606.99999999999, 22, -2.176078380049506, 48.13110481914438),
    @(45613.99999999999, 22, -1.698587494548043, 46.60073149377291),
    @(45620.99999999999, 23, -1.870587713303603, 48.67962470672012),
    @(45627.99999999999, 23, -4.470330134543165, 47.1717534601454),
    @(45634.99999999999, 24, -2.015736029566218, 48.5119160643932)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
}

$wsForecast.Range("A1").Select()
